$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume/change (E) figures to the refreshed
# values pulled by the GitHub Actions symbol-list update.
# Cells are stored as literal text (not numbers) in the source data, so we
# force Text number-format before writing the value (this preserves things
# like trailing zeros and the leading/trailing % sign exactly), then reset
# the cell style back to Normal so no stray formatting is introduced.
$updates = @{
    "D2" = "303.91"
    "D3" = "37.23"
    "E3" = "3.67%"
    "D4" = "5.037"
    "E4" = "-0.95%"
    "E5" = "-0.43%"
    "D6" = "2.211"
    "E6" = "-3.60%"
    "D7" = "7.996"
    "E7" = "-1.18%"
    "E8" = "0.58%"
    "D9" = "0.9287"
    "E9" = "0.28%"
    "D10" = "0.09842"
    "E10" = "-2.21%"
    "D11" = "0.1887"
    "E11" = "3.23%"
    "D12" = "0.08712"
    "E12" = "0.61%"
    "D13" = "0.03613"
    "E13" = "5.89%"
    "D14" = "0.09946"
    "E14" = "0.42%"
    "D15" = "0.001481"
    "E15" = "-0.08%"
    "D16" = "0.005648"
    "E16" = "0.91%"
    "D17" = "3.463"
    "E17" = "-0.81%"
    "D18" = "2.312"
    "E18" = "10.35%"
    "E20" = "0.60%"
    "D21" = "4.755"
    "E21" = "4.31%"
    "E22" = "-1.38%"
    "D23" = "0.04601"
    "E23" = "-1.30%"
    "D24" = "0.005204"
    "E24" = "15.85%"
    "D25" = "0.001253"
    "E25" = "1.02%"
    "D26" = "0.0001400"
    "E26" = "7.82%"
    "E39" = "4.21%"
    "D40" = "0.04764"
    "E40" = "1.37%"
    "D41" = "0.007929"
    "E41" = "1.09%"
    "E42" = "-1.17%"
    "D43" = "0.007564"
    "E43" = "-13.96%"
    "D44" = "0.002250"
    "E44" = "1.55%"
    "D45" = "0.01040"
    "E45" = "13.35%"
    "D46" = "0.00006306"
    "E46" = "4.88%"
    "D47" = "0.00000000750"
    "E47" = "0.13%"
    "D49" = "31.15"
    "E49" = "699.20%"
    "D51" = "0.00002100"
    "E51" = "0.13%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
